# Rename the resource item type used on the "Meta" sheet, and move the
# active selection down to B2 to reflect the edited cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Meta")

$ws.Range("B2").Value = "Hydra.Infrastructure.I18n.ResourceItem, Hydra.Infrastructure"

$ws.Activate()
$ws.Range("B2").Select()
